$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of the Barbaria online slot game and play for free. Enjoy stunning graphic features, two bonuses, and a medium volatility level.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml)

# ------------------------------------------------------------------
# 2) Near the bottom of the document there are two trailing paragraphs:
#      - a bold duplicate of the title ("Play Barbaria Free - ...")
#      - an italic meta-description paragraph ("Read our review of ...")
#    Replace both of them (in one shot, to dodge end-of-document range
#    quirks) with a single italic paragraph holding the new image-prompt
#    text.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldPara = $d.Paragraphs.Item($count - 1)
$italicPara = $d.Paragraphs.Item($count)
$tailRange = $d.Range($boldPara.Range.Start, $italicPara.Range.End)

$tailXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Create a feature image fitting the game &quot;Barbaria&quot; with the following specifications: - A cartoon-style image - Features a happy Maya warrior with glasses - The environment should depict an ancient Roman setting with barbarian weapons in the background. The image should be bright and colorful, with a focus on the Maya warrior as the main character. The warrior should look happy and adventurous, holding a sword or an axe in one hand and wearing glasses. The background can be of an ancient Roman setting with barbarian weapons displayed. The image should be able to attract players looking for a thrilling game with great graphics and visuals.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$tailRange.InsertXML($tailXml)
